# Applies the cryptos-list refresh described by the commit diff.
# Updates Price/Volume(1h) figures and re-orders a handful of coin rows
# (by overwriting Coin/Link/Price/Volume cells in place).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '52.045.89'
$ws.Range("E2").Value = '  -14.23%  '

# Row 3
$ws.Range("D3").Value = '2.272.10'
$ws.Range("E3").Value = '  -21.79%  '

# Row 5
$ws.Range("D5").Value = '''428.52'
$ws.Range("E5").Value = '  -18.87%  '

# Row 6
$ws.Range("D6").Value = '''118.07'
$ws.Range("E6").Value = '  -17.66%  '

# Row 7
$ws.Range("D7").Value = '''1.00'
$ws.Range("E7").Value = '  +0.29%  '

# Row 8
$ws.Range("D8").Value = '''0.464'
$ws.Range("E8").Value = '  -16.47%  '

# Row 9
$ws.Range("D9").Value = '2.271.85'
$ws.Range("E9").Value = '  -21.98%  '

# Row 10
$ws.Range("D10").Value = '''5.19'
$ws.Range("E10").Value = '  -14.00%  '

# Row 11
$ws.Range("D11").Value = '''0.0889'
$ws.Range("E11").Value = '  -17.65%  '

# Row 12
$ws.Range("D12").Value = '''0.300'
$ws.Range("E12").Value = '  -16.79%  '

# Row 13
$ws.Range("D13").Value = '''0.120'
$ws.Range("E13").Value = '  -5.28%  '

# Row 14
$ws.Range("D14").Value = '2.665.37'
$ws.Range("E14").Value = '  -22.06%  '

# Row 15
$ws.Range("D15").Value = '52.072.97'
$ws.Range("E15").Value = '  -14.14%  '

# Row 16
$ws.Range("D16").Value = '''18.75'
$ws.Range("E16").Value = '  -17.68%  '

# Row 17
$ws.Range("D17").Value = '''0.0000117'
$ws.Range("E17").Value = '  -17.10%  '

# Row 18
$ws.Range("D18").Value = '2.285.60'
$ws.Range("E18").Value = '  -21.62%  '

# Row 19
$ws.Range("D19").Value = '''3.91'
$ws.Range("E19").Value = '  -22.28%  '

# Row 20
$ws.Range("D20").Value = '''295.37'
$ws.Range("E20").Value = '  -18.41%  '

# Row 21
$ws.Range("D21").Value = '''8.83'
$ws.Range("E21").Value = '  -24.70%  '

# Row 22
$ws.Range("D22").Value = '''0.998'
$ws.Range("E22").Value = '  -0.21%  '

# Row 23
$ws.Range("E23").Value = '  -1.75%  '

# Row 24
$ws.Range("D24").Value = '''5.26'
$ws.Range("E24").Value = '  -20.80%  '

# Row 25
$ws.Range("D25").Value = '''54.38'
$ws.Range("E25").Value = '  -16.01%  '

# Row 26
$ws.Range("E26").Value = '  +0.30%  '

# Row 27
$ws.Range("D27").Value = '''0.147'
$ws.Range("E27").Value = '  -17.98%  '

# Row 28
$ws.Range("D28").Value = '''0.360'
$ws.Range("E28").Value = '  -20.86%  '

# Row 29
$ws.Range("B29").Value = 'USDe'
$ws.Range("C29").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D29").Value = '''0.999'
$ws.Range("E29").Value = '  -0.05%  '

# Row 30
$ws.Range("B30").Value = 'InternetComputer(DFINITY)'
$ws.Range("C30").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D30").Value = '''6.79'
$ws.Range("E30").Value = '  -13.82%  '

# Row 31
$ws.Range("D31").Value = '0.0₃0684'
$ws.Range("E31").Value = '  -19.47%  '

# Row 32
$ws.Range("D32").Value = '''141.45'
$ws.Range("E32").Value = '  -6.23%  '

# Row 33
$ws.Range("D33").Value = '''16.80'
$ws.Range("E33").Value = '  -15.12%  '

# Row 34
$ws.Range("D34").Value = '''1.31'
$ws.Range("E34").Value = '  -22.09%  '

# Row 35
$ws.Range("D35").Value = '''4.65'
$ws.Range("E35").Value = '  -16.74%  '

# Row 36
$ws.Range("D36").Value = '''3.46'
$ws.Range("E36").Value = '  -20.74%  '

# Row 37
$ws.Range("D37").Value = '''0.806'
$ws.Range("E37").Value = '  -19.45%  '

# Row 38
$ws.Range("B38").Value = 'FirstDigitalUSD'
$ws.Range("C38").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D38").Value = '''0.995'
$ws.Range("E38").Value = '  -0.16%  '

# Row 39
$ws.Range("B39").Value = 'ImmutableX'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D39").Value = '''0.983'
$ws.Range("E39").Value = '  -18.67%  '

# Row 40
$ws.Range("B40").Value = 'OKB'
$ws.Range("C40").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D40").Value = '''32.67'
$ws.Range("E40").Value = '  -13.71%  '

# Row 41
$ws.Range("D41").Value = '''10.25'
$ws.Range("E41").Value = '  -0.64%  '

# Row 42
$ws.Range("D42").Value = '''3.11'
$ws.Range("E42").Value = '  -16.30%  '

# Row 43
$ws.Range("B43").Value = 'Stacks'
$ws.Range("C43").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D43").Value = '''1.20'
$ws.Range("E43").Value = '  -19.47%  '

# Row 44
$ws.Range("B44").Value = 'Hedera'
$ws.Range("C44").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D44").Value = '''0.0490'
$ws.Range("E44").Value = '  -16.30%  '

# Row 45
$ws.Range("D45").Value = '1.860.35'
$ws.Range("E45").Value = '  -18.94%  '

# Row 46
$ws.Range("B46").Value = 'VeChain'
$ws.Range("C46").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D46").Value = '''0.0205'
$ws.Range("E46").Value = '  -13.71%  '

# Row 47
$ws.Range("B47").Value = 'Mantle'
$ws.Range("C47").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D47").Value = '''0.507'
$ws.Range("E47").Value = '  -21.85%  '

# Row 48
$ws.Range("D48").Value = '''0.0820'
$ws.Range("E48").Value = '  -11.16%  '

# Row 49
$ws.Range("D49").Value = '''15.31'
$ws.Range("E49").Value = '  -25.09%  '

# Row 51
$ws.Range("D51").Value = '''3.78'
$ws.Range("E51").Value = '  -24.72%  '

Write-Output "Edit applied successfully"
